$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 1 and row 2 (header row <-> id row), columns A:I
for ($col = 1; $col -le 9; $col++) {
    $cell1 = $ws.Cells.Item(1, $col)
    $cell2 = $ws.Cells.Item(2, $col)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
}

# Move the selection to A6 (as reflected in the saved file)
$ws.Range("A6").Select() | Out-Null
